# AP90_TestData_CreateandApplyOrUnapplyPrepayments_21C.xlsx
# "Add files via upload" - Anu scrubbed the live Oracle Cloud URL /
# credentials (UserName/Password) that had accidentally been left as a
# hyperlinked cell block on the Input_Value sheet, clearing Z2:AB2 and
# removing the hyperlink that pointed at the environment URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate() | Out-Null

# Remove the hyperlink anchored at Z2 (pointed at the Oracle Cloud env URL)
foreach ($h in @($ws.Hyperlinks)) {
    $h.Delete() | Out-Null
}

# Clear out the credential/URL cells that the hyperlink lived alongside;
# the cell formatting (fill/font/borders) is left intact, only the
# contained values are removed.
$ws.Range("Z2").Value = ""
$ws.Range("AA2").Value = ""
$ws.Range("AB2").Value = ""

# Leave the view focused on the now-empty block, matching where the
# author's selection ended up after the cleanup.
$ws.Range("Z2:AB2").Select() | Out-Null
